$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RF)
$ws.Range("B3").Value = 0.276
$ws.Range("C3").Value = -0.044
$ws.Range("D3").Value = 0.459
$ws.Range("E3").Value = 0.677
$ws.Range("F3").Value = 0.768
$ws.Range("G3").Value = 0.5649999999999999

# Row 4 (NN)
$ws.Range("B4").Value = 0.216
$ws.Range("C4").Value = -0.131
$ws.Range("D4").Value = 0.498
$ws.Range("E4").Value = 0.706
$ws.Range("F4").Value = 0.744
$ws.Range("G4").Value = 0.514

# Row 5 (RNN)
$ws.Range("B5").Value = -0.004
$ws.Range("C5").Value = -0.24
$ws.Range("D5").Value = 0.573
$ws.Range("E5").Value = 0.757
$ws.Range("F5").Value = 0.748
$ws.Range("G5").Value = 0.478
